$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.084852457046509
$ws.Range("B1").Value = 2.449600696563721
$ws.Range("C1").Value = 2.959489345550537
$ws.Range("D1").Value = 5.225256443023682
$ws.Range("E1").Value = 3.43882942199707
